# Update the cryptocurrency price/volume snapshot (GitHub Actions refresh).
# Most rows just get refreshed Price (column D) and Volume(1h) (column E)
# values. Rows 29/30 (Toncoin <-> PancakeSwap) and rows 40/41
# (LidoDAOToken <-> Celestia) also swap places in the ranking, so their
# Coin name (B) and Link (C) are updated too.
#
# Numeric-looking values in columns D/E are written with a leading
# apostrophe so Excel keeps them as text (matching the original
# inlineStr/text cells) instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.174.37"
$ws.Range("E2").Value = "'  -1.90%  "

$ws.Range("D3").Value = "2.141.95"
$ws.Range("E3").Value = "'  -3.37%  "

$ws.Range("E4").Value = "'  +0.01%  "

$ws.Range("D5").Value = "'235.02"
$ws.Range("E5").Value = "'  -2.80%  "

$ws.Range("E6").Value = "'  -4.89%  "

$ws.Range("D7").Value = "'69.22"
$ws.Range("E7").Value = "'  -5.72%  "

$ws.Range("E8").Value = "'  +0.04%  "

$ws.Range("D9").Value = "'0.565"
$ws.Range("E9").Value = "'  -7.01%  "

$ws.Range("D10").Value = "'38.47"
$ws.Range("E10").Value = "'  -10.15%  "

$ws.Range("D11").Value = "'0.0891"
$ws.Range("E11").Value = "'  -7.07%  "

$ws.Range("D12").Value = "'53.22"
$ws.Range("E12").Value = "'  -7.39%  "

$ws.Range("D13").Value = "'0.0993"
$ws.Range("E13").Value = "'  -4.17%  "

$ws.Range("D14").Value = "'6.57"
$ws.Range("E14").Value = "'  -6.75%  "

$ws.Range("D15").Value = "2.460.74"
$ws.Range("E15").Value = "'  -3.47%  "

$ws.Range("D16").Value = "'14.24"
$ws.Range("E16").Value = "'  -0.22%  "

$ws.Range("D17").Value = "2.125.43"
$ws.Range("E17").Value = "'  -3.47%  "

$ws.Range("D18").Value = "'0.774"
$ws.Range("E18").Value = "'  -7.66%  "

$ws.Range("D19").Value = "41.046.48"
$ws.Range("E19").Value = "'  -2.07%  "

$ws.Range("D20").Value = "0.0₃0990"
$ws.Range("E20").Value = "'  -7.65%  "

$ws.Range("D21").Value = "'68.76"
$ws.Range("E21").Value = "'  -5.37%  "

$ws.Range("D22").Value = "'5.70"
$ws.Range("E22").Value = "'  -7.71%  "

$ws.Range("D23").Value = "'224.38"
$ws.Range("E23").Value = "'  -2.58%  "

$ws.Range("D24").Value = "'9.41"
$ws.Range("E24").Value = "'  -12.55%  "

$ws.Range("E25").Value = "'  +0.09%  "

$ws.Range("D26").Value = "'1.89"
$ws.Range("E26").Value = "'  -8.67%  "

$ws.Range("D27").Value = "'10.50"
$ws.Range("E27").Value = "'  -10.79%  "

$ws.Range("E28").Value = "'  -9.58%  "

# Row 29 now holds PancakeSwap (was Toncoin)
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E29").Value = "'  -6.42%  "

# Row 30 now holds Toncoin (was PancakeSwap)
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.15"
$ws.Range("E30").Value = "'  -1.72%  "

$ws.Range("D31").Value = "'169.35"
$ws.Range("E31").Value = "'  +0.74%  "

$ws.Range("D32").Value = "'19.52"
$ws.Range("E32").Value = "'  -4.90%  "

$ws.Range("D33").Value = "'30.87"
$ws.Range("E33").Value = "'  +1.86%  "

$ws.Range("D34").Value = "'0.0748"
$ws.Range("E34").Value = "'  -6.27%  "

$ws.Range("D35").Value = "'5.02"
$ws.Range("E35").Value = "'  -11.89%  "

$ws.Range("E36").Value = "'  -5.19%  "

$ws.Range("D37").Value = "'0.101"
$ws.Range("E37").Value = "'  -8.48%  "

$ws.Range("D38").Value = "'4.08"
$ws.Range("E38").Value = "'  -4.50%  "

$ws.Range("D39").Value = "'0.0282"
$ws.Range("E39").Value = "'  -7.53%  "

# Row 40 now holds Celestia (was LidoDAOToken)
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'11.66"
$ws.Range("E40").Value = "'  -16.92%  "

# Row 41 now holds LidoDAOToken (was Celestia)
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "'2.03"
$ws.Range("E41").Value = "'  -4.38%  "

$ws.Range("D42").Value = "'5.21"
$ws.Range("E42").Value = "'  -7.73%  "

$ws.Range("D43").Value = "'57.25"
$ws.Range("E43").Value = "'  -12.97%  "

$ws.Range("D44").Value = "'0.185"
$ws.Range("E44").Value = "'  -7.00%  "

$ws.Range("D45").Value = "'8.16"
$ws.Range("E45").Value = "'  -8.19%  "

$ws.Range("D46").Value = "'0.0951"
$ws.Range("E46").Value = "'  -5.51%  "

$ws.Range("D47").Value = "'96.71"
$ws.Range("E47").Value = "'  -8.05%  "

$ws.Range("D48").Value = "'1.07"
$ws.Range("E48").Value = "'  -4.56%  "

$ws.Range("D49").Value = "'1.10"
$ws.Range("E49").Value = "'  -6.62%  "

$ws.Range("E50").Value = "'  -3.37%  "

$ws.Range("E51").Value = "'  -11.69%  "
